# Update the handoff/handback datetimes on the "zh-cn" and "de-de" sheets
# to reflect the regenerated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 17:08:48"
$wsZhCn.Range("H2").Value = "2016-03-13 17:09:05"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 17:08:52"
$wsDeDe.Range("H2").Value = "2016-03-13 17:09:11"
